$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab from "Cart Report" to "Report"
$ws.Name = "Report"

# Update the summary title text
$ws.Range("A1").Value = "Cart Report Summary"

# Turn row 3 into a header row for the Field/Value table
$ws.Range("A3").Value = "Field"
$ws.Range("B3").Value = "Value"

# Add the data row underneath the header
$ws.Range("A4").Value = "Number of products added to the cart"

# Store the count as text "1" (not a numeric 1) -- format the cell as
# Text first so Excel keeps the digit-only entry as a string value
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1"

# Set explicit column widths (best-fit / custom width) matching the target layout
$ws.Columns.Item(1).ColumnWidth = 35.11328125
$ws.Columns.Item(2).ColumnWidth = 6.09375
